$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new date row below the existing "Date" header.
# Force the cell to be treated as text (shared string) first, so that
# Excel does not auto-convert the date-like string into a numeric date value.
$cell = $ws.Range("A2")
$cell.NumberFormat = "@"
$cell.Value = "08.01.2025"
# Apply a date display format to the (still text-typed) cell, matching the source file.
$cell.NumberFormat = "d/mm/yyyy"

# Update the sheet selection to the newly added cell, as in the authored edit.
$cell.Select()
